$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Organic & Biomolecular Chemistry (row 38): add print and online ISSN
$ws.Range("C38").Value = "1477-0520"
$ws.Range("D38").Value = "1477-0539"

# Chemical Science (row 8): add online ISSN
$ws.Range("D8").Value = "2041-6539"

# CrystEngComm (row 11): add online ISSN
$ws.Range("D11").Value = "1466-8033"

[void]$ws.Range("D11").Select()
